# Update "想去人数" (want-to-go count) figures in the 展览 (Exhibition) and
# 全部类型 (All types) sheets to reflect a freshly regenerated data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 26
$ws1.Range("F7").Value  = 1043
$ws1.Range("F8").Value  = 519
$ws1.Range("F13").Value = 334
$ws1.Range("F15").Value = 319
$ws1.Range("F16").Value = 418
$ws1.Range("F17").Value = 5388
$ws1.Range("F19").Value = 1508
$ws1.Range("F20").Value = 342
$ws1.Range("F21").Value = 4480
$ws1.Range("F24").Value = 1431
$ws1.Range("F29").Value = 3770

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 26
$ws4.Range("F10").Value = 1043
$ws4.Range("F11").Value = 519
$ws4.Range("F16").Value = 334
$ws4.Range("F22").Value = 319
$ws4.Range("F24").Value = 418
$ws4.Range("F25").Value = 5388
$ws4.Range("F27").Value = 1508
$ws4.Range("F30").Value = 342
$ws4.Range("F32").Value = 4480
$ws4.Range("F35").Value = 1431
$ws4.Range("F45").Value = 3770
